$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3: was blank -> -14.2
$ws.Range("D3").Value = -14.2

# Update D5: was -14.4 -> blank
$ws.Range("D5").Value = $null

# Update D21: was blank -> -14.3
$ws.Range("D21").Value = -14.3

# Update D23: was -13.9 -> blank
$ws.Range("D23").Value = $null

# Delete entire row 26 (RM 232) - rows below shift up
$ws.Range("A26").EntireRow.Delete()

# After the first deletion, old row 28 (SC 92) is now row 27. Delete it too.
$ws.Range("A27").EntireRow.Delete()

# Now the row that was old row 34 (SC 193) is at row 32. Its D value was blank, set to -14.7
$ws.Range("D32").Value = -14.7
